# Apply the update described by the diff to the
# "Impacts_per_ingredients" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Shorten the header labels in row 1 (columns B..O)
# ---------------------------------------------------------------
$ws.Range("B1").Value = "Land Use Arable"
$ws.Range("C1").Value = "Land Use Fallow"
$ws.Range("D1").Value = "Land Use Perm Past"
$ws.Range("E1").Value = "GHG LUC"
$ws.Range("F1").Value = "GHG Feed"
$ws.Range("G1").Value = "GHG Farm"
$ws.Range("H1").Value = "GHG Processing"
$ws.Range("I1").Value = "GHG Transport"
$ws.Range("J1").Value = "GHG Packging"
$ws.Range("K1").Value = "GHG Retail"
$ws.Range("L1").Value = "Acidification"
$ws.Range("M1").Value = "Eutrophication"
$ws.Range("N1").Value = "Freshwater Withdrawals (FW)"
$ws.Range("O1").Value = "Scarcity-Weighted FW"

# ---------------------------------------------------------------
# 2. Recipe 1 (Beef + Rice, rows 2-3) - no more "Total" row,
#    recipe name renamed to "Rice and Beef"
# ---------------------------------------------------------------
$ws.Range("P2").Value = "Rice and Beef"
$ws.Range("P3").Value = "Rice and Beef"

# ---------------------------------------------------------------
# 3. Recipe 2 becomes "Rice and chicken" (Rice + Chicken, rows 4-5)
#    replacing the old Total row (4) / Rice row (5)
# ---------------------------------------------------------------
$ws.Range("A4").Value = "Rice"
$ws.Range("B4").Value = 0.4186
$ws.Range("C4").Value = 0.0616
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.004399999999999999
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.7106
$ws.Range("H4").Value = 0.013
$ws.Range("I4").Value = 0.0192
$ws.Range("J4").Value = 0.0168
$ws.Range("K4").Value = 0.0126
$ws.Range("L4").Value = 0.004714
$ws.Range("M4").Value = 0.005944
$ws.Range("N4").Value = 392.4
$ws.Range("O4").Value = 8655
$ws.Range("P4").Value = "Rice and chicken"

$ws.Range("A5").Value = "Chicken"
$ws.Range("B5").Value = 2.5753
$ws.Range("C5").Value = 0.83265
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1.1557
$ws.Range("F5").Value = 0.807625
$ws.Range("G5").Value = 0.30576
$ws.Range("H5").Value = 0.2002
$ws.Range("I5").Value = 0.12558
$ws.Range("J5").Value = 0.09645999999999999
$ws.Range("K5").Value = 0.080535
$ws.Range("L5").Value = 0.02929745
$ws.Range("M5").Value = 0.01345435
$ws.Range("N5").Value = 182.91
$ws.Range("O5").Value = 4016.74
$ws.Range("P5").Value = "Rice and chicken"

# ---------------------------------------------------------------
# 4. Recipe 3 becomes "Rice and Soy tofu" (Rice + Soybean Tofu,
#    rows 6-7) replacing the old Rice row (6) / Total row (7)
# ---------------------------------------------------------------
$ws.Range("A6").Value = "Rice"
$ws.Range("B6").Value = 0.4186
$ws.Range("C6").Value = 0.0616
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.004399999999999999
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.7106
$ws.Range("H6").Value = 0.013
$ws.Range("I6").Value = 0.0192
$ws.Range("J6").Value = 0.0168
$ws.Range("K6").Value = 0.0126
$ws.Range("L6").Value = 0.004714
$ws.Range("M6").Value = 0.005944
$ws.Range("N6").Value = 392.4
$ws.Range("O6").Value = 8655
$ws.Range("P6").Value = "Rice and Soy tofu"

$ws.Range("A7").Value = "Soybean Tofu"
$ws.Range("B7").Value = 1.125215
$ws.Range("C7").Value = 0.309855
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.43589
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.225225
$ws.Range("H7").Value = 0.36127
$ws.Range("I7").Value = 0.080535
$ws.Range("J7").Value = 0.080535
$ws.Range("K7").Value = 0.12285
$ws.Range("L7").Value = 0.00284375
$ws.Range("M7").Value = 0.00251615
$ws.Range("N7").Value = 63.245
$ws.Range("O7").Value = 2177.175
$ws.Range("P7").Value = "Rice and Soy tofu"

# ---------------------------------------------------------------
# 5. Remove the now-obsolete rows 8-10 (old Recipe 3 block),
#    shrinking the sheet from A1:P10 down to A1:P7
# ---------------------------------------------------------------
$ws.Range("A8:P10").EntireRow.Delete()
